$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of row number -> new sentiment value in column L
$updates = @{
    2  = "NEUTRAL"
    3  = "NEUTRAL"
    4  = "NEUTRAL"
    5  = "NEUTRAL"
    6  = "NEUTRAL"
    7  = "NEUTRAL"
    9  = "NEUTRAL"
    11 = "NEUTRAL"
    12 = "NEUTRAL"
    13 = "NEUTRAL"
    15 = "NEUTRAL"
    16 = "NEUTRAL"
    18 = "NEUTRAL"
    19 = "NEUTRAL"
    20 = "POSITIVE"
    21 = "POSITIVE"
}

foreach ($row in $updates.Keys) {
    $ws.Range("L$row").Value = $updates[$row]
}
